# Apply the changes described by the diff:
#  1. Update the cached "datetimeFigureOut" field text from 2020-10-13 to
#     2023-09-28 everywhere it appears (slide master + every slide layout's
#     Date placeholder).
#  2. Fix the typo in slide 7's title: "Reply Attacks" -> "Replay Attacks".

$p = $ppt.ActivePresentation

$oldDate = "2020-10-13"
$newDate = "2023-09-28"

function Update-DateShape($shapes) {
    for ($k = 1; $k -le $shapes.Count; $k++) {
        $shp = $shapes.Item($k)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide Master date placeholder.
$master = $p.SlideMaster
Update-DateShape $master.Shapes

# Every slide layout's date placeholder.
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    Update-DateShape $layout.Shapes
}

# Slide 7 title typo fix: "Reply Attacks" -> "Replay Attacks".
$slide7 = $p.Slides.Item(7)
$title = $slide7.Shapes.Item(1)
$titleRange = $title.TextFrame.TextRange
$titleRange.Text = $titleRange.Text.Replace("Reply Attacks", "Replay Attacks")
